# Regenerate save_data column G ("K", formerly "Strike#") values.
# For each data row (2-43) the new K value is written; row 24 is unchanged (already 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 5
    22 = 1
    23 = 2
    25 = 0
    26 = 1
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 0
    40 = 2
    41 = 0
    42 = 1
    43 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
